$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in B2:C4
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 4

$ws.Range("B3").Value = 3

$ws.Range("B4").Value = 1

# Add new rows 6 and 7
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1440
$ws.Range("C6").Value = 24

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1440
$ws.Range("C7").Value = 6

# Update selection to A8
$ws.Range("A8").Select()
